$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (it currently sits right after the
#    first paragraph's "：" run).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Move to the very end of the document (the trailing empty paragraph) and
#    append the copyright notice as three separate runs, matching the
#    formatting used in the source diff.
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$endRange.InsertAfter("Copyright ©2021-2099 ")
$endRange.Collapse(0)

$endRange.InsertAfter("HanxiaoZhang")
$endRange.Font.Name = "宋体"

$endRange.Collapse(0)
$endRange.InsertAfter(". All rights reserved")

# 3. Re-create the "_GoBack" bookmark as a zero-length bookmark right at the
#    new end of the document.
$endRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endRange)
